{"js": "// Update the two-digit multiplication answers in the table to a new\n// set of problems/answers, matching each old value to its replacement.\nconst replacements = [\n  [\"47\u00d751=2397\", \"39\u00d730=1170\"],\n  [\"31\u00d726=806\", \"15\u00d787=1305\"],\n  [\"52\u00d724=1248\", \"29\u00d759=1711\"],\n  [\"91\u00d747=4277\", \"48\u00d770=3360\"],\n  [\"17\u00d793=1581\", \"90\u00d774=6660\"],\n  [\"74\u00d711=814\", \"75\u00d726=1950\"],\n  [\"29\u00d736=1044\", \"17\u00d782=1394\"],\n  [\"50\u00d739=1950\", \"28\u00d722=616\"],\n  [\"20\u00d794=1880\", \"80\u00d747=3760\"],\n  [\"47\u00d771=3337\", \"35\u00d716=560\"],\n  [\"15\u00d715=225\", \"99\u00d733=3267\"],\n  [\"72\u00d735=2520\", \"94\u00d779=7426\"],\n  [\"51\u00d787=4437\", \"99\u00d712=1188\"],\n  [\"60\u00d762=3720\", \"52\u00d797=5044\"],\n  [\"26\u00d720=520\", \"13\u00d768=884\"],\n  [\"48\u00d727=1296\", \"32\u00d735=1120\"],\n  [\"67\u00d732=2144\", \"14\u00d762=868\"],\n  [\"18\u00d732=576\", \"91\u00d790=8190\"],\n  [\"66\u00d730=1980\", \"11\u00d734=374\"],\n  [\"59\u00d741=2419\", \"22\u00d717=374\"],\n  [\"85\u00d783=7055\", \"87\u00d772=6264\"],\n  [\"86\u00d727=2322\", \"11\u00d731=341\"],\n  [\"70\u00d714=980\", \"69\u00d795=6555\"],\n  [\"62\u00d744=2728\", \"57\u00d756=3192\"],\n  [\"83\u00d733=2739\", \"99\u00d729=2871\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication answers in the table to a new\n# set of problems/answers, matching each old value to its replacement.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"47\u00d751=2397\"; New = \"39\u00d730=1170\" },\n    @{ Old = \"31\u00d726=806\";  New = \"15\u00d787=1305\" },\n    @{ Old = \"52\u00d724=1248\"; New = \"29\u00d759=1711\" },\n    @{ Old = \"91\u00d747=4277\"; New = \"48\u00d770=3360\" },\n    @{ Old = \"17\u00d793=1581\"; New = \"90\u00d774=6660\" },\n    @{ Old = \"74\u00d711=814\";  New = \"75\u00d726=1950\" },\n    @{ Old = \"29\u00d736=1044\"; New = \"17\u00d782=1394\" },\n    @{ Old = \"50\u00d739=1950\"; New = \"28\u00d722=616\" },\n    @{ Old = \"20\u00d794=1880\"; New = \"80\u00d747=3760\" },\n    @{ Old = \"47\u00d771=3337\"; New = \"35\u00d716=560\" },\n    @{ Old = \"15\u00d715=225\";  New = \"99\u00d733=3267\" },\n    @{ Old = \"72\u00d735=2520\"; New = \"94\u00d779=7426\" },\n    @{ Old = \"51\u00d787=4437\"; New = \"99\u00d712=1188\" },\n    @{ Old = \"60\u00d762=3720\"; New = \"52\u00d797=5044\" },\n    @{ Old = \"26\u00d720=520\";  New = \"13\u00d768=884\" },\n    @{ Old = \"48\u00d727=1296\"; New = \"32\u00d735=1120\" },\n    @{ Old = \"67\u00d732=2144\"; New = \"14\u00d762=868\" },\n    @{ Old = \"18\u00d732=576\";  New = \"91\u00d790=8190\" },\n    @{ Old = \"66\u00d730=1980\"; New = \"11\u00d734=374\" },\n    @{ Old = \"59\u00d741=2419\"; New = \"22\u00d717=374\" },\n    @{ Old = \"85\u00d783=7055\"; New = \"87\u00d772=6264\" },\n    @{ Old = \"86\u00d727=2322\"; New = \"11\u00d731=341\" },\n    @{ Old = \"70\u00d714=980\";  New = \"69\u00d795=6555\" },\n    @{ Old = \"62\u00d744=2728\"; New = \"57\u00d756=3192\" },\n    @{ Old = \"83\u00d733=2739\"; New = \"99\u00d729=2871\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
